# Update header labels to shorter column names: "Album Title" -> "Album",
# "Song Title" -> "Song". Commit message: "Added concatenated song data"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Album"
$ws.Range("B1").Value = "Song"
